$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 50 data for Tahj Brooks ---
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "Tahj Brooks"
$ws.Range("C50").Value = "RB"
$ws.Range("D50").Value = 20
$ws.Range("E50").Value = 14
$ws.Range("F50").Value = 46
$ws.Range("G50").Value = "Benchwarmer"
$ws.Range("H50").Value = "Trayveon Williams, Ray Davis, Jeremy McNichols"
$ws.Range("I50").Value = "5'9"
$ws.Range("J50").Value = 214
$ws.Range("K50").Value = " Tahj Brooks is built to see a lot of work. He was only 5’9’’ at the combine, but he weighed `r`nin at 214 pounds, giving him a top-three Body Mass Index in this year’s class.`r`n Brooks spent five seasons at Texas Tech, and he was really featured in the offense during `r`nhis final two. In 2023 and 2024, Brooks saw 290 and 284 rush attempts, respectively, and he `r`nhit total yards per team play rates of 1.71 and 1.99. That 1.99 number ranks top-10 in the `r`nclass.`r`n There may be some long speed concerns with Brooks, even if he exceeded expectations a bit `r`nat the combine with a 4.52 40. He was below average at generating explosive runs throughout `r`nhis career. His receiving profile isn’t stellar, either, with a bottom-10 best-season `r`nreceiving yards per team pass attempt rate. Brooks’ 35.7 Breakout Score is seventh-worst in `r`nthe class.`r`n We’ve had successes in the ZAP Model database from backs who had sub-40 Breakout Scores, `r`nbut those wins haven’t typically been long-term ones. Here’s a list of backs with a `r`nBreakout Score below 40 to have given fantasy managers 14 or more PPR points per game in `r`none of their first three seasons in the league: Alfred Morris, Andre Ellington, Carlos `r`nHyde, Melvin Gordon, Jordan Howard, Chris Carson, Miles Sanders, and James Robinson. There `r`nare a lot of one-hit wonders in there. The best back in the group, Melvin Gordon, at least `r`nhad draft capital backing him, too.`r`n That same filter -- backs with Breakout Scores below 40 -- also hasn’t provided a lot of `r`nreceiving upside. Just nine backs in the model with that low of a Breakout Score ended up `r`nreaching a 10% target share per game rate in one of their first three seasons. `r`nBrooks can handle a lot of work -- he’s not a bad prospect -- but I’m not sure he’ll get `r`nthe opportunity to be a true three-down back in the NFL. "

# Notes column on this row wraps text, like the rest of the table
$ws.Range("K50").WrapText = $true

# Row grew tall enough to hit Excel's max row height, matching the other
# long scouting reports in this sheet
$ws.Rows.Item(50).RowHeight = 409.5

# Leave the selection where the editor ended up after typing the new row
$ws.Range("F47").Select()
